$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header ---
$ws.Range("B1").Value = "Values"
# Columns H1:L1 are no longer part of the table — remove them entirely
$ws.Range("H1:L2").Clear()

# --- Row 2 (existing row, data changes) ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Final Value"
$ws.Range("C2").Value = 2645071.686008946
$ws.Range("D2").Value = 3049452.098302247
$ws.Range("E2").Value = 2469548.226738013
$ws.Range("F2").Value = 2653360.104686474
$ws.Range("G2").Value = 2653245.833758925

# --- Row 3 (new row) — copy A2's format for the A-column style ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Annualized Return"
$ws.Range("C3").Value = 0.3788614333539579
$ws.Range("D3").Value = 0.4451950695318403
$ws.Range("E3").Value = 0.3479439652969061
$ws.Range("F3").Value = 0.3802869612158055
$ws.Range("G3").Value = 0.3802673280104469

# --- Row 4 (new row) — copy A2's format for the A-column style ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Sharpe Ratio"
$ws.Range("C4").Value = 0.730558553332663
$ws.Range("D4").Value = 0.9740564215860639
$ws.Range("E4").Value = 0.7671793834188858
$ws.Range("F4").Value = 0.7279514498061527
$ws.Range("G4").Value = 0.727610642345788

$excel.CutCopyMode = 0
